$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.423951666666667
$ws.Range("H2").Value = 16.271855
$ws.Range("I2").Value = 0.4774188439413272
$ws.Range("J2").Value = 0.4774188439413271
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 14.80997966666667
$ws.Range("N2").Value = 44.429939
$ws.Range("O2").Value = 0.2388798507865045
$ws.Range("P2").Value = 0.2388798507865045
$ws.Range("Q2").Value = 80.32861389631611
$ws.Range("R2").Value = 722.9575250668451
$ws.Range("S2").Value = 0.1140457422033697
$ws.Range("T2").Value = 0.1140457422033697
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.423951666666667
$ws.Range("H3").Value = 16.271855
$ws.Range("I3").Value = 0.4774188439413272
$ws.Range("J3").Value = 0.4774188439413271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.42883866666667
$ws.Range("N3").Value = 100.286516
$ws.Range("O3").Value = 0.5391956081231261
$ws.Range("P3").Value = 0.5391956081231262
$ws.Range("Q3").Value = 181.3164052007978
$ws.Range("R3").Value = 1631.84764680718
$ws.Range("S3").Value = 0.2574221438883837
$ws.Range("T3").Value = 0.2574221438883837
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.423951666666667
$ws.Range("H4").Value = 16.271855
$ws.Range("I4").Value = 0.4774188439413272
$ws.Range("J4").Value = 0.4774188439413271
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.509909333333333
$ws.Range("N4").Value = 7.529728
$ws.Range("O4").Value = 0.04048396962919451
$ws.Range("P4").Value = 0.04048396962919452
$ws.Range("Q4").Value = 13.61362691171555
$ws.Range("R4").Value = 122.52264220544
$ws.Range("S4").Value = 0.01932780997852585
$ws.Range("T4").Value = 0.01932780997852585
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.423951666666667
$ws.Range("H5").Value = 16.271855
$ws.Range("I5").Value = 0.4774188439413272
$ws.Range("J5").Value = 0.4774188439413271
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 11.24888166666667
$ws.Range("N5").Value = 33.746645
$ws.Range("O5").Value = 0.1814405714611748
$ws.Range("P5").Value = 0.1814405714611748
$ws.Range("Q5").Value = 61.01339046405278
$ws.Range("R5").Value = 549.120514176475
$ws.Range("S5").Value = 0.08662314787104784
$ws.Range("T5").Value = 0.08662314787104784
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.583504333333333
$ws.Range("H6").Value = 4.750513
$ws.Range("I6").Value = 0.1393808158066948
$ws.Range("J6").Value = 0.1393808158066948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.80997966666667
$ws.Range("N6").Value = 44.429939
$ws.Range("O6").Value = 0.2388798507865045
$ws.Range("P6").Value = 0.2388798507865045
$ws.Range("Q6").Value = 23.45166697874522
$ws.Range("R6").Value = 211.065002808707
$ws.Range("S6").Value = 0.03329526848240453
$ws.Range("T6").Value = 0.03329526848240452
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.583504333333333
$ws.Range("H7").Value = 4.750513
$ws.Range("I7").Value = 0.1393808158066948
$ws.Range("J7").Value = 0.1393808158066948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.42883866666667
$ws.Range("N7").Value = 100.286516
$ws.Range("O7").Value = 0.5391956081231261
$ws.Range("P7").Value = 0.5391956081231262
$ws.Range("Q7").Value = 52.93471088696756
$ws.Range("R7").Value = 476.412397982708
$ws.Range("S7").Value = 0.07515352373958822
$ws.Range("T7").Value = 0.07515352373958824
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.583504333333333
$ws.Range("H8").Value = 4.750513
$ws.Range("I8").Value = 0.1393808158066948
$ws.Range("J8").Value = 0.1393808158066948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.509909333333333
$ws.Range("N8").Value = 7.529728
$ws.Range("O8").Value = 0.04048396962919451
$ws.Range("P8").Value = 0.04048396962919452
$ws.Range("Q8").Value = 3.974452305607111
$ws.Range("R8").Value = 35.770070750464
$ws.Range("S8").Value = 0.005642688714010587
$ws.Range("T8").Value = 0.005642688714010587
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.583504333333333
$ws.Range("H9").Value = 4.750513
$ws.Range("I9").Value = 0.1393808158066948
$ws.Range("J9").Value = 0.1393808158066948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 11.24888166666667
$ws.Range("N9").Value = 33.746645
$ws.Range("O9").Value = 0.1814405714611748
$ws.Range("P9").Value = 0.1814405714611748
$ws.Range("Q9").Value = 17.81265286432055
$ws.Range("R9").Value = 160.313875778885
$ws.Range("S9").Value = 0.02528933487069145
$ws.Range("T9").Value = 0.02528933487069145
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6660723333333333
$ws.Range("H10").Value = 1.998217
$ws.Range("I10").Value = 0.05862800830537802
$ws.Range("J10").Value = 0.05862800830537802
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.80997966666667
$ws.Range("N10").Value = 44.429939
$ws.Range("O10").Value = 0.2388798507865045
$ws.Range("P10").Value = 0.2388798507865045
$ws.Range("Q10").Value = 9.864517713195889
$ws.Range("R10").Value = 88.78065941876301
$ws.Range("S10").Value = 0.01400504987589865
$ws.Range("T10").Value = 0.01400504987589865
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6660723333333333
$ws.Range("H11").Value = 1.998217
$ws.Range("I11").Value = 0.05862800830537802
$ws.Range("J11").Value = 0.05862800830537802
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 33.42883866666667
$ws.Range("N11").Value = 100.286516
$ws.Range("O11").Value = 0.5391956081231261
$ws.Range("P11").Value = 0.5391956081231262
$ws.Range("Q11").Value = 22.26602457133022
$ws.Range("R11").Value = 200.394221141972
$ws.Range("S11").Value = 0.03161196459126599
$ws.Range("T11").Value = 0.03161196459126599
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6660723333333333
$ws.Range("H12").Value = 1.998217
$ws.Range("I12").Value = 0.05862800830537802
$ws.Range("J12").Value = 0.05862800830537802
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.509909333333333
$ws.Range("N12").Value = 7.529728
$ws.Range("O12").Value = 0.04048396962919451
$ws.Range("P12").Value = 0.04048396962919452
$ws.Range("Q12").Value = 1.671781166108444
$ws.Range("R12").Value = 15.046030494976
$ws.Range("S12").Value = 0.002373494507655088
$ws.Range("T12").Value = 0.002373494507655088
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6660723333333333
$ws.Range("H13").Value = 1.998217
$ws.Range("I13").Value = 0.05862800830537802
$ws.Range("J13").Value = 0.05862800830537802
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 11.24888166666667
$ws.Range("N13").Value = 33.746645
$ws.Range("O13").Value = 0.1814405714611748
$ws.Range("P13").Value = 0.1814405714611748
$ws.Range("Q13").Value = 7.492568859107222
$ws.Range("R13").Value = 67.43311973196499
$ws.Range("S13").Value = 0.01063749933055829
$ws.Range("T13").Value = 0.01063749933055829
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.687463666666666
$ws.Range("H14").Value = 11.062391
$ws.Range("I14").Value = 0.3245723319466
$ws.Range("J14").Value = 0.3245723319466
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 14.80997966666667
$ws.Range("N14").Value = 44.429939
$ws.Range("O14").Value = 0.2388798507865045
$ws.Range("P14").Value = 0.2388798507865045
$ws.Range("Q14").Value = 54.61126192490545
$ws.Range("R14").Value = 491.5013573241491
$ws.Range("S14").Value = 0.0775337902248316
$ws.Range("T14").Value = 0.07753379022483162
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.687463666666666
$ws.Range("H15").Value = 11.062391
$ws.Range("I15").Value = 0.3245723319466
$ws.Range("J15").Value = 0.3245723319466
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 33.42883866666667
$ws.Range("N15").Value = 100.286516
$ws.Range("O15").Value = 0.5391956081231261
$ws.Range("P15").Value = 0.5391956081231262
$ws.Range("Q15").Value = 123.2676280021951
$ws.Range("R15").Value = 1109.408652019756
$ws.Range("S15").Value = 0.1750079759038881
$ws.Range("T15").Value = 0.1750079759038881
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.687463666666666
$ws.Range("H16").Value = 11.062391
$ws.Range("I16").Value = 0.3245723319466
$ws.Range("J16").Value = 0.3245723319466
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.509909333333333
$ws.Range("N16").Value = 7.529728
$ws.Range("O16").Value = 0.04048396962919451
$ws.Range("P16").Value = 0.04048396962919452
$ws.Range("Q16").Value = 9.255199473294221
$ws.Range("R16").Value = 83.29679525964801
$ws.Range("S16").Value = 0.01313997642900299
$ws.Range("T16").Value = 0.01313997642900299
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.687463666666666
$ws.Range("H17").Value = 11.062391
$ws.Range("I17").Value = 0.3245723319466
$ws.Range("J17").Value = 0.3245723319466
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 11.24888166666667
$ws.Range("N17").Value = 33.746645
$ws.Range("O17").Value = 0.1814405714611748
$ws.Range("P17").Value = 0.1814405714611748
$ws.Range("Q17").Value = 41.47984243646611
$ws.Range("R17").Value = 373.318581928195
$ws.Range("S17").Value = 0.05889058938887722
$ws.Range("T17").Value = 0.05889058938887723
